# "timelog, Video28, Indexist kõik asjad viidus MeasurePage'i"
# Fill in the rest of row 3 & 4 (sheet rows 6-8) of the "Nädal 7" timesheet:
#  - row 6 (entry 3): add end time, delta minutes, and rename the activity
#    comment from "V28" to "V28,29"
#  - row 7 (entry 4): add date/start/end time, delta minutes, switch the
#    activity to "apps", add a comment "trello, slack" and mark column C (t="s" v17 -> "x")
#  - row 8 (entry 5): add start time and an interruption time of 10 minutes
# The Total Time formula (F17) recalculates automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Nädal 7")

# Widen column H (Comments) so the new longer comment text is readable.
$ws.Columns.Item(8).ColumnWidth = 21.109375

# Row 6 (entry #3) — add the stop time & delta, update the comment text.
$ws.Range("D6").Value = 0.64583333333333337
$ws.Range("F6").Value = 75
$ws.Range("H6").Value = "V28,29"

# Row 7 (entry #4) — new full entry.
$ws.Range("B7").Value = 43909
$ws.Range("C7").Value = 0.64583333333333337
$ws.Range("D7").Value = 0.66666666666666663
$ws.Range("F7").Value = 30
$ws.Range("G7").Value = "apps"
$ws.Range("H7").Value = "trello, slack"
$ws.Range("I7").Value = "x"

# Row 8 (entry #5) — start time plus an interruption.
$ws.Range("C8").Value = 0.66666666666666663
$ws.Range("E8").Value = 10

# Leave the selection where the user ended up editing.
$ws.Range("E9").Select() | Out-Null
